$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4820.857
$ws.Range("I62").Value = 4223.375
$ws.Range("J62").Value = 5617.5
$ws.Range("K62").Value = 4223.375
$ws.Range("L62").Value = 5617.5
$ws.Range("M62").Value = -3599.375
$ws.Range("N62").Value = -6865.5

$ws.Range("H65").Value = 4820.857
$ws.Range("I65").Value = 4223.375
$ws.Range("J65").Value = 5617.5
$ws.Range("K65").Value = 21116.875
$ws.Range("L65").Value = 28087.5
$ws.Range("M65").Value = -17996.875
$ws.Range("N65").Value = -34327.5

$ws.Range("H129").Value = 772.3333
$ws.Range("I129").Value = 658.5
$ws.Range("J129").Value = 1000
$ws.Range("K129").Value = 1975.5
$ws.Range("L129").Value = 3000
$ws.Range("M129").Value = 3024.5
$ws.Range("N129").Value = -13000

$ws.Range("H135").Value = 637.10345
$ws.Range("I135").Value = 493
$ws.Range("K135").Value = 4437
$ws.Range("M135").Value = -1902

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2065.63
$ws.Range("I32").Value = 1913.4631
$ws.Range("K32").Value = 1913.4631
$ws.Range("M32").Value = -1626.4631

$ws.Range("H61").Value = 1245.0938
$ws.Range("I61").Value = 1077.9
$ws.Range("J61").Value = 1842.2142
$ws.Range("K61").Value = 1077.9
$ws.Range("L61").Value = 1842.2142
$ws.Range("M61").Value = -865.9000000000001
$ws.Range("N61").Value = -2266.2142

$ws.Range("H74").Value = 4728.125
$ws.Range("I74").Value = 5410.25
$ws.Range("J74").Value = 2681.75
$ws.Range("K74").Value = 5410.25
$ws.Range("L74").Value = 2681.75
$ws.Range("M74").Value = -4536.25
$ws.Range("N74").Value = -4429.75

$ws.Range("H77").Value = 4728.125
$ws.Range("I77").Value = 5410.25
$ws.Range("J77").Value = 2681.75
$ws.Range("K77").Value = 27051.25
$ws.Range("L77").Value = 13408.75
$ws.Range("M77").Value = -22683.25
$ws.Range("N77").Value = -22144.75

$ws.Range("H136").Value = 1245.0938
$ws.Range("I136").Value = 1077.9
$ws.Range("J136").Value = 1842.2142
$ws.Range("K136").Value = 3233.7
$ws.Range("L136").Value = 5526.642599999999
$ws.Range("M136").Value = -683.7000000000003
$ws.Range("N136").Value = -10626.6426

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1578.78
$ws.Range("I134").Value = 1356.1555
$ws.Range("J134").Value = 3582.4
$ws.Range("K134").Value = 4068.4665
$ws.Range("L134").Value = 10747.2
$ws.Range("M134").Value = -1533.4665
$ws.Range("N134").Value = -15817.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 34571.145
$ws.Range("J20").Value = 34571.145
$ws.Range("L20").Value = 34571.145
$ws.Range("N20").Value = -35043.145

$ws.Range("H30").Value = 34571.145
$ws.Range("J30").Value = 34571.145
$ws.Range("L30").Value = 34571.145
$ws.Range("N30").Value = -34753.145

$ws.Range("H99").Value = 3036.8635
$ws.Range("I99").Value = 3426.375
$ws.Range("J99").Value = 2814.2856
$ws.Range("K99").Value = 3426.375
$ws.Range("L99").Value = 2814.2856
$ws.Range("M99").Value = -1928.375
$ws.Range("N99").Value = -5810.2856

$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800

$ws.Range("H124").Value = 16250
$ws.Range("J124").Value = 16250
$ws.Range("L124").Value = 16250
$ws.Range("N124").Value = -21160

$ws.Range("H125").Value = 12993.333
$ws.Range("J125").Value = 12993.333
$ws.Range("L125").Value = 12993.333
$ws.Range("N125").Value = -17913.333

$ws.Range("H126").Value = 3036.8635
$ws.Range("I126").Value = 3426.375
$ws.Range("J126").Value = 2814.2856
$ws.Range("K126").Value = 10279.125
$ws.Range("L126").Value = 8442.856800000001
$ws.Range("M126").Value = -7809.125
$ws.Range("N126").Value = -13382.8568

$ws.Range("H128").Value = 34571.145
$ws.Range("J128").Value = 34571.145
$ws.Range("L128").Value = 34571.145
$ws.Range("N128").Value = -44531.145

$ws.Range("H130").Value = 26500
$ws.Range("J130").Value = 26500
$ws.Range("L130").Value = 26500
$ws.Range("N130").Value = -36540

$ws.Range("H131").Value = 29400
$ws.Range("J131").Value = 29400
$ws.Range("L131").Value = 29400
$ws.Range("N131").Value = -39480

$ws.Range("H134").Value = 1758.9375
$ws.Range("I134").Value = 1128.2245
$ws.Range("J134").Value = 3819.2666
$ws.Range("K134").Value = 3384.6735
$ws.Range("L134").Value = 11457.7998
$ws.Range("M134").Value = -849.6734999999999
$ws.Range("N134").Value = -16527.7998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 16070653
$ws.Range("I132").Value = 1450
$ws.Range("J132").Value = 19641586
$ws.Range("K132").Value = 13050
$ws.Range("L132").Value = 176774274
$ws.Range("M132").Value = -10520
$ws.Range("N132").Value = -176779334

$ws.Range("H138").Value = 1789.7826
$ws.Range("I138").Value = 725.5714
$ws.Range("K138").Value = 2176.7142
$ws.Range("M138").Value = 2963.2858

$ws.Range("H139").Value = 3408.9473
$ws.Range("I139").Value = 1193
$ws.Range("K139").Value = 3579
$ws.Range("M139").Value = 1561

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 180000
$ws.Range("J141").Value = 180000
$ws.Range("L141").Value = 180000
$ws.Range("N141").Value = -190360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1470.93
$ws.Range("I132").Value = 1498.8191
$ws.Range("J132").Value = 1034
$ws.Range("K132").Value = 4496.4573
$ws.Range("L132").Value = 3102
$ws.Range("M132").Value = -1966.4573
$ws.Range("N132").Value = -8162

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 30030
$ws.Range("J46").Value = 30030
$ws.Range("L46").Value = 30030
$ws.Range("N46").Value = -30492

$ws.Range("H134").Value = 30030
$ws.Range("J134").Value = 30030
$ws.Range("L134").Value = 90090
$ws.Range("N134").Value = -95160

$ws.Range("H136").Value = 1935.2258
$ws.Range("I136").Value = 1937.1945
$ws.Range("J136").Value = 1928.4762
$ws.Range("K136").Value = 5811.583500000001
$ws.Range("L136").Value = 5785.4286
$ws.Range("M136").Value = -3261.583500000001
$ws.Range("N136").Value = -10885.4286

$ws.Range("H141").Value = 145000
$ws.Range("J141").Value = 145000
$ws.Range("L141").Value = 145000
$ws.Range("N141").Value = -155360
